$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H4").Value = 410.0435
$ws.Range("I4").Value = 243.07692
$ws.Range("J4").Value = 627.1
$ws.Range("K4").Value = 243.07692
$ws.Range("L4").Value = 627.1
$ws.Range("M4").Value = -129.07692
$ws.Range("N4").Value = -855.1
$ws.Range("H112").Value = 927.3077
$ws.Range("I112").Value = 628.3333
$ws.Range("J112").Value = 1017
$ws.Range("K112").Value = 1884.9999
$ws.Range("L112").Value = 3051
$ws.Range("M112").Value = -776.9999
$ws.Range("N112").Value = -5267
$ws.Range("H113").Value = 2701.3635
$ws.Range("I113").Value = 1678.5714
$ws.Range("J113").Value = 3178.6667
$ws.Range("K113").Value = 1678.5714
$ws.Range("L113").Value = 3178.6667
$ws.Range("M113").Value = 1575.4286
$ws.Range("N113").Value = -9686.6667
$ws.Range("H116").Value = 2139660
$ws.Range("I116").Value = 2691011.2
$ws.Range("K116").Value = 2691011.2
$ws.Range("M116").Value = -2687569.2
$ws.Range("H127").Value = 844.76
$ws.Range("I127").Value = 521
$ws.Range("J127").Value = 1420.3334
$ws.Range("K127").Value = 1563
$ws.Range("L127").Value = 4261.0002
$ws.Range("M127").Value = 3397
$ws.Range("N127").Value = -14181.0002
$ws.Range("H129").Value = 6105.381
$ws.Range("I129").Value = 377.5
$ws.Range("J129").Value = 7895.3438
$ws.Range("K129").Value = 1132.5
$ws.Range("L129").Value = 23686.0314
$ws.Range("M129").Value = 3867.5
$ws.Range("N129").Value = -33686.0314
$ws.Range("H132").Value = 3287.8115
$ws.Range("I132").Value = 3094.9648
$ws.Range("K132").Value = 9284.894400000001
$ws.Range("M132").Value = -6754.894400000001
$ws.Range("H138").Value = 3222.9314
$ws.Range("I138").Value = 1750.9474
$ws.Range("J138").Value = 3740.8518
$ws.Range("K138").Value = 5252.8422
$ws.Range("L138").Value = 11222.5554
$ws.Range("M138").Value = -112.8422
$ws.Range("N138").Value = -21502.5554

$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 198043.92
$ws.Range("I74").Value = 2144.4707
$ws.Range("J74").Value = 589842.8
$ws.Range("K74").Value = 2144.4707
$ws.Range("L74").Value = 589842.8
$ws.Range("M74").Value = -1270.4707
$ws.Range("N74").Value = -591590.8
$ws.Range("H77").Value = 198043.92
$ws.Range("I77").Value = 2144.4707
$ws.Range("J77").Value = 589842.8
$ws.Range("K77").Value = 10722.3535
$ws.Range("L77").Value = 2949214
$ws.Range("M77").Value = -6354.353499999999
$ws.Range("N77").Value = -2957950
$ws.Range("H102").Value = 1763.4286
$ws.Range("J102").Value = 1799.5
$ws.Range("L102").Value = 1799.5
$ws.Range("N102").Value = -5043.5

$ws = $wb.Worksheets.Item(3)
$ws.Range("H94").Value = 15509.048
$ws.Range("I94").Value = 10994.8
$ws.Range("J94").Value = 19612.908
$ws.Range("K94").Value = 10994.8
$ws.Range("L94").Value = 19612.908
$ws.Range("M94").Value = -10543.8
$ws.Range("N94").Value = -20514.908
$ws.Range("H99").Value = 2334.7334
$ws.Range("J99").Value = 2391.3333
$ws.Range("L99").Value = 2391.3333
$ws.Range("N99").Value = -5387.3333
$ws.Range("H105").Value = 2218.75
$ws.Range("I105").Value = 1912.5
$ws.Range("K105").Value = 1912.5
$ws.Range("M105").Value = -165.5

$ws = $wb.Worksheets.Item(4)
$ws.Range("H132").Value = 2140.9546
$ws.Range("I132").Value = 1728.5
$ws.Range("J132").Value = 3997
$ws.Range("K132").Value = 5185.5
$ws.Range("L132").Value = 11991
$ws.Range("M132").Value = -2655.5
$ws.Range("N132").Value = -17051
$ws.Range("H134").Value = 4633.433
$ws.Range("I134").Value = 5059.864
$ws.Range("J134").Value = 3460.75
$ws.Range("K134").Value = 15179.592
$ws.Range("L134").Value = 10382.25
$ws.Range("M134").Value = -12644.592
$ws.Range("N134").Value = -15452.25

$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 330041.66
$ws.Range("I2").Value = 495012.5
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 2970075
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -2969962
$ws.Range("N2").Value = -826
$ws.Range("H5").Value = 781.7708
$ws.Range("I5").Value = 467.04166
$ws.Range("J5").Value = 1096.5
$ws.Range("K5").Value = 1401.12498
$ws.Range("L5").Value = 3289.5
$ws.Range("M5").Value = -1289.12498
$ws.Range("N5").Value = -3513.5
$ws.Range("H68").Value = 1716
$ws.Range("I68").Value = 400
$ws.Range("K68").Value = 1200
$ws.Range("M68").Value = -389
$ws.Range("H71").Value = 1716
$ws.Range("I71").Value = 400
$ws.Range("K71").Value = 3600
$ws.Range("M71").Value = 456
$ws.Range("H113").Value = 1063.4407
$ws.Range("I113").Value = 472
$ws.Range("J113").Value = 1867.8
$ws.Range("K113").Value = 1416
$ws.Range("L113").Value = 5603.4
$ws.Range("M113").Value = 754
$ws.Range("N113").Value = -9943.4
$ws.Range("H131").Value = 1755106.5
$ws.Range("J131").Value = 979.9545000000001
$ws.Range("L131").Value = 2939.8635
$ws.Range("N131").Value = -13019.8635
$ws.Range("H134").Value = 893.6539
$ws.Range("I134").Value = 639.6875
$ws.Range("K134").Value = 1919.0625
$ws.Range("M134").Value = 3150.9375
$ws.Range("H135").Value = 781.7708
$ws.Range("I135").Value = 467.04166
$ws.Range("J135").Value = 1096.5
$ws.Range("K135").Value = 4203.37494
$ws.Range("L135").Value = 9868.5
$ws.Range("M135").Value = -1668.37494
$ws.Range("N135").Value = -14938.5
$ws.Range("H136").Value = 1203.6875
$ws.Range("H137").Value = 4443
$ws.Range("I137").Value = 786.61536
$ws.Range("J137").Value = 12365.167
$ws.Range("K137").Value = 2359.84608
$ws.Range("L137").Value = 37095.501
$ws.Range("M137").Value = 2740.15392
$ws.Range("N137").Value = -47295.501
$ws.Range("H140").Value = 1397
$ws.Range("I140").Value = 942.2308
$ws.Range("K140").Value = 2826.6924
$ws.Range("M140").Value = 2353.3076

$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 1266.5454
$ws.Range("I113").Value = 1191.625
$ws.Range("J113").Value = 1466.3334
$ws.Range("K113").Value = 1191.625
$ws.Range("L113").Value = 1466.3334
$ws.Range("M113").Value = 978.375
$ws.Range("N113").Value = -5806.3334
$ws.Range("H126").Value = 2781.2307
$ws.Range("I126").Value = 1826.375
$ws.Range("J126").Value = 4309
$ws.Range("K126").Value = 5479.125
$ws.Range("L126").Value = 12927
$ws.Range("M126").Value = -3009.125
$ws.Range("N126").Value = -17867

$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 2564799.8
$ws.Range("I22").Value = 4762103
$ws.Range("J22").Value = 1279.5
$ws.Range("K22").Value = 4762103
$ws.Range("L22").Value = 1279.5
$ws.Range("M22").Value = -4761808
$ws.Range("N22").Value = -1869.5
$ws.Range("H27").Value = 2564799.8
$ws.Range("I27").Value = 4762103
$ws.Range("J27").Value = 1279.5
$ws.Range("K27").Value = 4762103
$ws.Range("L27").Value = 1279.5
$ws.Range("M27").Value = -4761996
$ws.Range("N27").Value = -1493.5
$ws.Range("H93").Value = 2428.4285
$ws.Range("J93").Value = 1999.5
$ws.Range("L93").Value = 1999.5
$ws.Range("N93").Value = -4495.5
$ws.Range("H122").Value = 4270
$ws.Range("I122").Value = 4195.3125
$ws.Range("K122").Value = 12585.9375
$ws.Range("M122").Value = -10135.9375
$ws.Range("H136").Value = 8675.75
$ws.Range("I136").Value = 11772.167
$ws.Range("J136").Value = 4031.125
$ws.Range("K136").Value = 35316.501
$ws.Range("L136").Value = 12093.375
$ws.Range("M136").Value = -32766.501
$ws.Range("N136").Value = -17193.375

$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 1126.3549
$ws.Range("I122").Value = 948.34784
$ws.Range("J122").Value = 1638.125
$ws.Range("K122").Value = 2845.04352
$ws.Range("L122").Value = 4914.375
$ws.Range("M122").Value = -395.0435200000002
$ws.Range("N122").Value = -9814.375
